# ActivationFlow.pptx - slide 1 layout tweaks (positions/sizes, connector
# flips, and OnX/ActivationService label font sizes 8pt -> 16pt).
#
# NOTE on literal point values below: Shape.Left/Top/Width/Height are
# expressed in points in the COM object model, while the underlying OOXML
# stores EMU (1 pt = 12700 EMU). The runtime's points<->EMU round trip loses
# precision (float32-ish truncation), so naive `emu/12700` literals can miss
# the exact target EMU by 1 unit. The literals used here were solved so that
# they land on the exact target EMU after the runtime's internal conversion.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id=4 "OnLaunched" rounded rectangle: shift left ---
$sh4 = $s.Shapes.Item(1)
$sh4.Left = 57.86567115783692
$tr4 = $sh4.TextFrame.TextRange
$tr4.Font.Size = 16

# --- Shape id=5 "Straight Arrow Connector 4": flip horizontally + reroute ---
$sh5 = $s.Shapes.Item(2)
$sh5.HorizontalFlip = -1
$sh5.Left = 279.49998474121094
$sh5.Top = 458.3220367431641
$sh5.Width = 0.18559055775403976
$sh5.Height = 71.54921340942384

# --- Shape id=6 "OnActivated" rounded rectangle: shift + widen ---
$sh6 = $s.Shapes.Item(3)
$sh6.Left = 216.6194534301758
$sh6.Width = 126.13220596313478
$tr6 = $sh6.TextFrame.TextRange
$tr6.Font.Size = 16

# --- Shape id=7 "OnBackgroundActivated" rounded rectangle: shift + widen ---
$sh7 = $s.Shapes.Item(4)
$sh7.Left = 387.55110168457037
$sh7.Width = 126.13220596313478
$tr7 = $sh7.TextFrame.TextRange
$tr7.Font.Size = 16

# --- Shape id=8 "ActivationService" / "ActivateAsync" rounded rectangle ---
$sh8 = $s.Shapes.Item(5)
$sh8.Left = 190.99999237060547
$sh8.Width = 176.99991607666018
$tr8 = $sh8.TextFrame.TextRange
$tr8.Font.Size = 16

# --- Shape id=9 "Straight Connector 8": shift left ---
$sh9 = $s.Shapes.Item(6)
$sh9.Left = 114.84267807006836

# --- Shape id=10 "Straight Connector 9": un-flip + reroute ---
$sh10 = $s.Shapes.Item(7)
$sh10.HorizontalFlip = 0
$sh10.Left = 450.6172332763672
$sh10.Top = 459.4585113525391
$sh10.Width = 0
$sh10.Height = 35.77440834045411

# --- Shape id=11 "Straight Connector 10": shift left + lengthen ---
$sh11 = $s.Shapes.Item(8)
$sh11.Left = 114.84267807006836
$sh11.Width = 335.77449035644537
